# devops-getting-started.pptx
# "Change customer to stakeholder and all people to people."
#
# Slide 2 ("Values") has two text edits:
#   1. TextBox 58 - "It's about delighting our customers with VALUE!"
#      -> "customers" becomes "stakeholders".
#   2. TextBox 67 - "All PEOPLE need to buy into the transformation..."
#      -> the leading "All " run is dropped so the sentence reads
#         "PEOPLE need to buy into the transformation...".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- 1) "customers" -> "stakeholders" -------------------------------------
# (The run's text is rewritten in full, using the literal right single
#  quotation mark U+2019 used by the original "It's" - the TextRange
#  getter normalizes curly quotes to straight ones, so round-tripping the
#  existing text through Get+Replace would corrupt that character.)
$apostrophe = [char]0x2019
$valueShape = $s.Shapes.Item("TextBox 58")
$introRun = $valueShape.TextFrame.TextRange.Paragraphs(1).Runs(1)
$introRun.Text = "It" + $apostrophe + "s about delighting our stakeholders with "

# --- 2) "All PEOPLE ..." -> "PEOPLE ..." -----------------------------------
$peopleShape = $s.Shapes.Item("TextBox 67")
$peopleParagraph = $peopleShape.TextFrame.TextRange.Paragraphs(3)
$leadingAllRun = $peopleParagraph.Runs(1)
if ($leadingAllRun.Text -eq "All ") {
    $leadingAllRun.Text = ""
}
